$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "through" date references for 2022-08-03
$ws.Name = "Through 2022-08-03"
$ws.Range("I1").Value = "2022 (through 08-03)"

# Update the monthly data points that changed
$ws.Range("I7").Value = 142
$ws.Range("I9").Value = 16
$ws.Range("I14").Value = 986
